$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.927.40"
$ws.Range("E2").Value = "  +5.89%  "
$ws.Range("D3").Value = "2.976.23"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$orig = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.33"
$ws.Range("D5").Style = $orig
$ws.Range("E5").Value = "  +2.42%  "
$orig = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.04"
$ws.Range("D6").Style = $orig
$ws.Range("E6").Value = "  +7.08%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "2.971.79"
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E10").Value = "  +3.92%  "
$ws.Range("E11").Value = "  +2.77%  "
$orig = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("D12").Style = $orig
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("E13").Value = "  +1.62%  "
$orig = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.90"
$ws.Range("D14").Style = $orig
$ws.Range("E14").Value = "  +6.32%  "
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "64.921.31"
$ws.Range("E16").Value = "  +5.74%  "
$ws.Range("D17").Value = "3.470.90"
$ws.Range("E17").Value = "  +3.07%  "
$orig = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.88"
$ws.Range("D18").Style = $orig
$ws.Range("E18").Value = "  +3.93%  "
$ws.Range("D19").Value = "2.992.73"
$ws.Range("E19").Value = "  +3.65%  "
$orig = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "447.95"
$ws.Range("D20").Style = $orig
$ws.Range("E20").Value = "  +3.61%  "
$orig = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("D21").Style = $orig
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("E23").Value = "  +5.27%  "
$orig = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.86"
$ws.Range("D24").Style = $orig
$ws.Range("E24").Value = "  +2.01%  "
$orig = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.25"
$ws.Range("D25").Style = $orig
$ws.Range("E25").Value = "  +4.07%  "
$orig = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.57"
$ws.Range("D26").Style = $orig
$ws.Range("E26").Value = "  +5.92%  "
$orig = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.18"
$ws.Range("D27").Style = $orig
$ws.Range("E27").Value = "  +8.34%  "
$ws.Range("E28").Value = "  -0.03%  "
$orig = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.82"
$ws.Range("D29").Style = $orig
$ws.Range("E29").Value = "  +11.93%  "
$orig = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.36"
$ws.Range("D30").Style = $orig
$ws.Range("E30").Value = "  +14.77%  "
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("E32").Value = "  -1.01%  "
$orig = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("D33").Style = $orig
$ws.Range("E33").Value = "  +3.57%  "
$orig = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.54"
$ws.Range("D34").Style = $orig
$ws.Range("E34").Value = "  +3.60%  "
$orig = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = $orig
$orig = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.980"
$ws.Range("D36").Style = $orig
$ws.Range("E36").Value = "  +2.66%  "
$ws.Range("E37").Value = "  +4.18%  "
$orig = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.08"
$ws.Range("D38").Style = $orig
$ws.Range("E38").Value = "  +7.93%  "
$orig = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.90"
$ws.Range("D39").Style = $orig
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  +1.83%  "
$orig = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "43.74"
$ws.Range("D41").Style = $orig
$ws.Range("E41").Value = "  +10.46%  "
$ws.Range("E42").Value = "  +5.23%  "
$orig = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.297"
$ws.Range("D43").Style = $orig
$ws.Range("E43").Value = "  +11.20%  "
$orig = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.39"
$ws.Range("D44").Style = $orig
$ws.Range("E44").Value = "  +1.77%  "
$orig = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "379.40"
$ws.Range("D45").Style = $orig
$ws.Range("E45").Value = "  +11.91%  "
$ws.Range("D46").Value = "2.760.40"
$ws.Range("E46").Value = "  +2.82%  "
$orig = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0347"
$ws.Range("D47").Style = $orig
$ws.Range("E47").Value = "  +4.10%  "
$orig = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.89"
$ws.Range("D48").Style = $orig
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +2.19%  "
$orig = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.95"
$ws.Range("D51").Style = $orig
$ws.Range("E51").Value = "  +7.05%  "
